$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.ClearFormats()
}

Set-TextValue "D2" "28.288.77"
Set-TextValue "E2" "  +0.31%  "
Set-TextValue "D3" "1.858.48"
Set-TextValue "E3" "  -0.52%  "
Set-TextValue "D4" "0.9999"
Set-TextValue "E4" "  -0.01%  "
Set-TextValue "D5" "329.97"
Set-TextValue "E5" "  -2.12%  "
Set-TextValue "E6" "  +0.07%  "
Set-TextValue "D7" "0.4554"
Set-TextValue "E7" "  -3.34%  "
Set-TextValue "D8" "0.3932"
Set-TextValue "E8" "  +0.11%  "
Set-TextValue "D9" "47.48"
Set-TextValue "E9" "  +1.15%  "
Set-TextValue "D10" "0.07813"
Set-TextValue "E10" "  -1.99%  "
Set-TextValue "D11" "0.9819"
Set-TextValue "E11" "  -2.55%  "
Set-TextValue "D12" "21.27"
Set-TextValue "E12" "  -2.06%  "
Set-TextValue "D13" "1.850.11"
Set-TextValue "E13" "  -0.88%  "
Set-TextValue "D14" "5.807"
Set-TextValue "E14" "  -2.96%  "
Set-TextValue "D15" "6.953"
Set-TextValue "E15" "  -4.40%  "
Set-TextValue "D16" "1.000"
Set-TextValue "E16" "  -0.08%  "
Set-TextValue "D17" "88.04"
Set-TextValue "E17" "  -3.57%  "
Set-TextValue "D18" "0.06529"
Set-TextValue "D19" "0.00001017"
Set-TextValue "E19" "  -2.49%  "
Set-TextValue "D20" "17.05"
Set-TextValue "E20" "  -3.97%  "
Set-TextValue "D22" "28.267.43"
Set-TextValue "E22" "  +0.23%  "
Set-TextValue "D23" "5.297"
Set-TextValue "E23" "  -2.61%  "
Set-TextValue "D24" "10.74"
Set-TextValue "E24" "  -2.98%  "
Set-TextValue "D25" "2.250"
Set-TextValue "E25" "  -1.67%  "
Set-TextValue "D26" "2.086.21"
Set-TextValue "E26" "  +0.10%  "
Set-TextValue "D27" "156.51"
Set-TextValue "E27" "  -1.59%  "
Set-TextValue "D28" "19.19"
Set-TextValue "E28" "  -3.28%  "
Set-TextValue "D29" "2.050"
Set-TextValue "E29" "  -3.87%  "
Set-TextValue "D30" "5.268"
Set-TextValue "E30" "  -3.90%  "
Set-TextValue "D31" "116.21"
Set-TextValue "E31" "  -3.07%  "
Set-TextValue "B32" "Stellar"
Set-TextValue "C32" "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D32" "0.09270"
Set-TextValue "E32" "  -2.41%  "
Set-TextValue "B33" "ImmutableX"
Set-TextValue "C33" "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
Set-TextValue "D33" "0.9351"
Set-TextValue "E33" "  -4.02%  "
Set-TextValue "D34" "3.598"
Set-TextValue "E34" "  +0.68%  "
Set-TextValue "D35" "1.377"
Set-TextValue "E35" "  -0.07%  "
Set-TextValue "D36" "5.186"
Set-TextValue "E36" "  -2.99%  "
Set-TextValue "D37" "0.06001"
Set-TextValue "E37" "  -1.58%  "
Set-TextValue "D38" "0.02195"
Set-TextValue "E38" "  -3.42%  "
Set-TextValue "D39" "8.171"
Set-TextValue "E39" "  -2.91%  "
Set-TextValue "D40" "1.158"
Set-TextValue "E40" "  -1.59%  "
Set-TextValue "D41" "1.000"
Set-TextValue "E41" "  +0.08%  "
Set-TextValue "D42" "0.5672"
Set-TextValue "E42" "  -5.17%  "
Set-TextValue "D43" "9.968"
Set-TextValue "E43" "  -3.88%  "
Set-TextValue "D44" "0.1789"
Set-TextValue "E44" "  -5.03%  "
Set-TextValue "D45" "1.243"
Set-TextValue "E45" "  -2.46%  "
Set-TextValue "D46" "2.271"
Set-TextValue "E46" "  +18.18%  "
Set-TextValue "B47" "Decentraland"
Set-TextValue "C47" "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
Set-TextValue "D47" "0.5373"
Set-TextValue "E47" "  -4.50%  "
Set-TextValue "B48" "EnergySwap"
Set-TextValue "C48" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D48" "11.78"
Set-TextValue "E48" "  -3.27%  "
Set-TextValue "D49" "0.07161"
Set-TextValue "E49" "  +4.32%  "
Set-TextValue "D50" "1.859"
Set-TextValue "E50" "  -6.25%  "
Set-TextValue "D51" "109.56"
Set-TextValue "E51" "  -1.87%  "
